$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value2 = '27.716.19'
$ws.Cells.Item(2,5).Value2 = '  -1.07%  '
$ws.Cells.Item(3,4).Value2 = '1.883.66'
$ws.Cells.Item(3,5).Value2 = '  -0.36%  '
$ws.Cells.Item(4,4).Value2 = "'1.002"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Value2 = '  -0.02%  '
$ws.Cells.Item(5,4).Value2 = "'331.15"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value2 = '  +1.40%  '
$ws.Cells.Item(6,4).Value2 = "'1.001"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value2 = '  -0.08%  '
$ws.Cells.Item(7,4).Value2 = "'0.4748"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value2 = '  +3.62%  '
$ws.Cells.Item(8,4).Value2 = "'0.3986"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value2 = '  +1.04%  '
$ws.Cells.Item(9,4).Value2 = "'47.92"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value2 = '  -6.65%  '
$ws.Cells.Item(10,4).Value2 = "'0.08068"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value2 = '  -1.85%  '
$ws.Cells.Item(11,4).Value2 = "'1.026"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value2 = '  -1.22%  '
$ws.Cells.Item(12,4).Value2 = "'21.82"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value2 = '  +0.62%  '
$ws.Cells.Item(13,4).Value2 = '1.898.46'
$ws.Cells.Item(13,5).Value2 = '  +0.32%  '
$ws.Cells.Item(14,4).Value2 = "'5.977"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value2 = '  -0.16%  '
$ws.Cells.Item(15,4).Value2 = "'7.203"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value2 = '  -1.76%  '
$ws.Cells.Item(16,4).Value2 = "'1.002"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value2 = '  +0.02%  '
$ws.Cells.Item(17,4).Value2 = "'87.06"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value2 = '  -2.89%  '
$ws.Cells.Item(18,4).Value2 = "'0.00001041"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value2 = '  -1.35%  '
$ws.Cells.Item(19,4).Value2 = "'0.06608"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value2 = '  +0.40%  '
$ws.Cells.Item(20,4).Value2 = "'17.31"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value2 = '  -1.58%  '
$ws.Cells.Item(21,4).Value2 = "'1.001"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value2 = '  -0.13%  '
$ws.Cells.Item(22,4).Value2 = '27.730.90'
$ws.Cells.Item(22,5).Value2 = '  -1.00%  '
$ws.Cells.Item(23,4).Value2 = "'5.520"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value2 = '  -2.42%  '
$ws.Cells.Item(24,4).Value2 = "'11.00"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value2 = '  -1.01%  '
$ws.Cells.Item(25,4).Value2 = "'2.311"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value2 = '  +0.12%  '
$ws.Cells.Item(26,4).Value2 = '2.119.64'
$ws.Cells.Item(26,5).Value2 = '  -0.03%  '
$ws.Cells.Item(27,4).Value2 = "'155.42"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value2 = '  +0.81%  '
$ws.Cells.Item(28,4).Value2 = "'20.22"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value2 = '  +1.52%  '
$ws.Cells.Item(29,4).Value2 = "'2.099"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value2 = '  -0.59%  '
$ws.Cells.Item(30,4).Value2 = "'5.588"
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value2 = '  -1.77%  '
$ws.Cells.Item(31,4).Value2 = "'122.47"
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value2 = '  -1.50%  '
$ws.Cells.Item(32,4).Value2 = "'0.9702"
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value2 = '  +1.13%  '
$ws.Cells.Item(33,4).Value2 = "'0.09554"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value2 = '  +0.13%  '
$ws.Cells.Item(34,4).Value2 = "'1.468"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value2 = '  -0.46%  '
$ws.Cells.Item(35,4).Value2 = "'3.624"
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value2 = '  -0.21%  '
$ws.Cells.Item(36,4).Value2 = "'5.308"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value2 = '  -2.99%  '
$ws.Cells.Item(37,4).Value2 = "'0.06123"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value2 = '  +0.18%  '
$ws.Cells.Item(38,4).Value2 = "'0.02259"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value2 = '  -1.01%  '
$ws.Cells.Item(39,4).Value2 = "'1.230"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value2 = '  -1.71%  '
$ws.Cells.Item(40,4).Value2 = "'8.161"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value2 = '  -5.67%  '
$ws.Cells.Item(41,4).Value2 = "'0.6008"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value2 = '  -1.59%  '
$ws.Cells.Item(42,4).Value2 = "'1.001"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value2 = '  -0.16%  '
$ws.Cells.Item(43,4).Value2 = "'0.1898"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value2 = '  +0.42%  '
$ws.Cells.Item(44,4).Value2 = "'10.33"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value2 = '  -4.02%  '
$ws.Cells.Item(45,4).Value2 = "'1.253"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value2 = '  -4.18%  '
$ws.Cells.Item(46,4).Value2 = "'0.5700"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value2 = '  -1.94%  '
$ws.Cells.Item(47,4).Value2 = "'12.20"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value2 = '  -3.76%  '
$ws.Cells.Item(48,4).Value2 = "'3.408"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value2 = '  -0.50%  '
$ws.Cells.Item(49,4).Value2 = "'1.942"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value2 = '  -2.67%  '
$ws.Cells.Item(50,4).Value2 = "'0.06822"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value2 = '  -0.91%  '
$ws.Cells.Item(51,4).Value2 = "'110.81"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value2 = '  +0.48%  '
